$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append a new daily entry (2025-02-04) as rows 135-140, following the exact
# same layout/formatting pattern used by the previous day's block
# (rows 124-132: blank separator, date/Domm, Meeting/Reconsile,
#  General Discussion, Study/<topic>, Total).
#
# NOTE on ordering: values/formulas are written to each new cell BEFORE its
# formatting is copied over (PasteSpecial formats-only) from the matching
# template cell. Doing it the other way round (format copy, then later
# writing a value to a cell that a not-yet-evaluated formula elsewhere
# depends on) leaves that formula's cached result stale in this engine.
# ---------------------------------------------------------------------------

# --- Row 136: date row ---
$ws.Range("A136").Value = 45692
$ws.Range("B136").Value = "Domm"
$ws.Range("D136").Value = 0.25

# --- Row 137: Meeting / Reconsile row ---
$ws.Range("B137").Value = "Meeting"
$ws.Range("C137").Value = "Reconsile"
$ws.Range("D137").Value = 0

# --- Row 138: General Discussion row ---
$ws.Range("C138").Value = "General Discussion"
$ws.Range("D138").Value = 0.25

# --- Row 139: Study / ASP.NET Core Final Demo row ---
$ws.Range("B139").Value = "Study"
$ws.Range("C139").Value = "ASP.NET Core Final Demo"
$ws.Range("D139").Value = 7.5

# --- Row 140: Total row ---
$ws.Range("B140").Value = "Total"
$ws.Range("D140").Formula = "=SUM(D135:D139)"

# --- Now copy the formatting over from the prior day's equivalent block ---

# Row 135: blank separator row (format copied from row 124)
$ws.Range("A124:D124").Copy()
$ws.Range("A135:D135").PasteSpecial(-4122)   # xlPasteFormats

# Row 136 format (from row 125)
$ws.Range("A125:D125").Copy()
$ws.Range("A136:D136").PasteSpecial(-4122)

# Row 137 format (from row 126)
$ws.Range("A126:D126").Copy()
$ws.Range("A137:D137").PasteSpecial(-4122)

# Row 138 format (from row 127)
$ws.Range("A127:D127").Copy()
$ws.Range("A138:D138").PasteSpecial(-4122)

# Row 139 format (from row 128)
$ws.Range("A128:D128").Copy()
$ws.Range("A139:D139").PasteSpecial(-4122)

# Row 140 format (from row 132, A/B + D only -- no C cell in the new row)
$ws.Range("A132:B132").Copy()
$ws.Range("A140:B140").PasteSpecial(-4122)
$ws.Range("D132").Copy()
$ws.Range("D140").PasteSpecial(-4122)

# Clear clipboard marching-ants / leftover copy-mode reference
$excel.CutCopyMode = 0

# Reflect the scrolled/selected state after the edit (best-effort: the
# sheet was scrolled down and the newly added block selected).
$ws.Range("A136:D141").Select()
